# repull data, push all data, mean calculation
# Update column F (dSF) values on the active sheet to reflect the
# repulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = 1
    4  = -5
    5  = 4
    6  = -3
    7  = -4
    8  = 1
    9  = 2
    12 = -1
    13 = 4
    14 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
